$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 59, shifting existing rows 59..80 down to 60..81
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record
$ws.Cells.Item(59, 1).Value = 2
$ws.Cells.Item(59, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44636
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112030
$ws.Cells.Item(59, 7).Value = "Poroto granado"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 600
$ws.Cells.Item(59, 11).Value = 18000
$ws.Cells.Item(59, 12).Value = 20000
$ws.Cells.Item(59, 13).Value = 19000
$ws.Cells.Item(59, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(59, 16).Value = 760
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# Match the date cell format used by the rest of column D
$ws.Cells.Item(59, 4).NumberFormat = $ws.Cells.Item(60, 4).NumberFormat
